$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The EC (account statement) periods are being reversed in order - most
# recent period first. Swap the period label (column E) and the
# "Valor Mora" amount (column F) between row 16 and row 23, 17 and 22,
# 18 and 21, 19 and 20 so that the whole block of rows is reversed.

$pairs = @(
    @(16, 23),
    @(17, 22),
    @(18, 21),
    @(19, 20)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $e1 = $ws.Cells.Item($r1, 5).Value2
    $e2 = $ws.Cells.Item($r2, 5).Value2
    $f1 = $ws.Cells.Item($r1, 6).Value2
    $f2 = $ws.Cells.Item($r2, 6).Value2

    $ws.Cells.Item($r1, 5).Value2 = $e2
    $ws.Cells.Item($r2, 5).Value2 = $e1
    $ws.Cells.Item($r1, 6).Value2 = $f2
    $ws.Cells.Item($r2, 6).Value2 = $f1
}
